# The document's hidden "_GoBack" bookmark (the marker Word drops at the
# location of the most recent edit) moved from the empty paragraph just
# before the "RequestInspection(MyKit mk){" method body, up to the empty
# paragraph right after "List<Arm> Arms //Only 4" (just before the
# "//Shared data with KitRobot..." paragraph). Re-create that by adding a
# fresh "_GoBack" bookmark at the new location -- Word only ever keeps a
# single "_GoBack" bookmark, so adding it here automatically removes the
# bookmark from its old spot.

$d = $word.ActiveDocument

$count = $d.Paragraphs.Count

# Locate the empty paragraph that follows "List<Arm> Arms //Only 4" --
# this is where the "_GoBack" bookmark needs to end up.
$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $text = $d.Paragraphs.Item($i).Range.Text
    if ($text -eq "List<Arm> Arms //Only 4`r") {
        $targetIndex = $i + 1
        break
    }
}

if ($targetIndex -ne -1) {
    $targetParagraph = $d.Paragraphs.Item($targetIndex)
    $d.Bookmarks.Add("_GoBack", $targetParagraph.Range)
}
